$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.152264333333333
$ws.Range("H2").Value = 3.456793
$ws.Range("I2").Value = 0.3930660006090215
$ws.Range("J2").Value = 0.3930660006090216
$ws.Range("M2").Value = 1.923239
$ws.Range("N2").Value = 5.769717
$ws.Range("O2").Value = 0.2340262838603868
$ws.Range("P2").Value = 0.2340262838603868
$ws.Range("Q2").Value = 2.216079704175667
$ws.Range("R2").Value = 19.944717337581
$ws.Range("S2").Value = 0.09198777543439385
$ws.Range("T2").Value = 0.09198777543439385
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.152264333333333
$ws.Range("H3").Value = 3.456793
$ws.Range("I3").Value = 0.3930660006090215
$ws.Range("J3").Value = 0.3930660006090216
$ws.Range("O3").Value = 0.4335574295612247
$ws.Range("P3").Value = 0.4335574295612246
$ws.Range("Q3").Value = 4.105512442433111
$ws.Range("R3").Value = 36.949611981898
$ws.Range("S3").Value = 0.1704166848719582
$ws.Range("T3").Value = 0.1704166848719582
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.152264333333333
$ws.Range("H4").Value = 3.456793
$ws.Range("I4").Value = 0.3930660006090215
$ws.Range("J4").Value = 0.3930660006090216
$ws.Range("O4").Value = 0.3324162865783886
$ws.Range("P4").Value = 0.3324162865783886
$ws.Range("Q4").Value = 3.147770301148222
$ws.Range("R4").Value = 28.329932710334
$ws.Range("S4").Value = 0.1306615403026696
$ws.Range("T4").Value = 0.1306615403026696
$ws.Range("I5").Value = 0.4761983545501621
$ws.Range("J5").Value = 0.476198354550162
$ws.Range("M5").Value = 1.923239
$ws.Range("N5").Value = 5.769717
$ws.Range("O5").Value = 0.2340262838603868
$ws.Range("P5").Value = 0.2340262838603868
$ws.Range("Q5").Value = 2.684774330635
$ws.Range("R5").Value = 24.162968975715
$ws.Range("S5").Value = 0.1114429312958053
$ws.Range("T5").Value = 0.1114429312958053
$ws.Range("I6").Value = 0.4761983545501621
$ws.Range("J6").Value = 0.476198354550162
$ws.Range("O6").Value = 0.4335574295612247
$ws.Range("P6").Value = 0.4335574295612246
$ws.Range("S6").Value = 0.206459334560053
$ws.Range("T6").Value = 0.2064593345600529
$ws.Range("I7").Value = 0.4761983545501621
$ws.Range("J7").Value = 0.476198354550162
$ws.Range("O7").Value = 0.3324162865783886
$ws.Range("P7").Value = 0.3324162865783886
$ws.Range("S7").Value = 0.1582960886943038
$ws.Range("T7").Value = 0.1582960886943038
$ws.Range("I8").Value = 0.1307356448408163
$ws.Range("J8").Value = 0.1307356448408163
$ws.Range("M8").Value = 1.923239
$ws.Range("N8").Value = 5.769717
$ws.Range("O8").Value = 0.2340262838603868
$ws.Range("P8").Value = 0.2340262838603868
$ws.Range("Q8").Value = 0.7370787824313332
$ws.Range("R8").Value = 6.633709041882
$ws.Range("S8").Value = 0.0305955771301876
$ws.Range("T8").Value = 0.0305955771301876
$ws.Range("I9").Value = 0.1307356448408163
$ws.Range("J9").Value = 0.1307356448408163
$ws.Range("O9").Value = 0.4335574295612247
$ws.Range("P9").Value = 0.4335574295612246
$ws.Range("S9").Value = 0.05668141012921352
$ws.Range("T9").Value = 0.05668141012921351
$ws.Range("I10").Value = 0.1307356448408163
$ws.Range("J10").Value = 0.1307356448408163
$ws.Range("O10").Value = 0.3324162865783886
$ws.Range("P10").Value = 0.3324162865783886
$ws.Range("S10").Value = 0.04345865758141524
$ws.Range("T10").Value = 0.04345865758141523
